$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => (D new value or $null, E new value or $null)
$updates = @(
    @{Row=2;  D="27.347.40";  E="  -0.71%  "},
    @{Row=3;  D="1.712.05";   E="  -0.64%  "},
    @{Row=4;  D="1.005";      E="  +0.01%  "},
    @{Row=5;  D="224.46";     E="  -0.73%  "},
    @{Row=6;  D="0.5289";     E="  -1.30%  "},
    @{Row=7;  D=$null;        E="  -0.01%  "},
    @{Row=8;  D="0.06684";    E="  +1.14%  "},
    @{Row=9;  D="0.2666";     E="  -0.15%  "},
    @{Row=10; D="20.89";      E="  -3.94%  "},
    @{Row=11; D="0.07688";    E="  -0.72%  "},
    @{Row=12; D="4.503";      E="  -2.71%  "},
    @{Row=13; D="1.946.82";   E="  -0.73%  "},
    @{Row=14; D="1.708.21";   E="  -0.83%  "},
    @{Row=15; D="0.5837";     E="  -0.25%  "},
    @{Row=16; D="0.0₅8219";   E="  -1.30%  "},
    @{Row=17; D="68.01";      E="  +0.04%  "},
    @{Row=18; D="27.361.75";  E="  -0.71%  "},
    @{Row=19; D="222.36";     E="  +0.87%  "},
    @{Row=20; D="1.005";      E="  +0.06%  "},
    @{Row=21; D="4.630";      E="  -2.17%  "},
    @{Row=22; D="10.42";      E="  -2.07%  "},
    @{Row=23; D=$null;        E="  -1.30%  "},
    @{Row=24; D=$null;        E="  -0.03%  "},
    @{Row=25; D="144.85";     E="  -2.70%  "},
    @{Row=26; D="1.687";      E="  -2.35%  "},
    @{Row=27; D="0.1206";     E="  -2.47%  "},
    @{Row=28; D="7.240";      E="  -2.39%  "},
    @{Row=29; D="16.27";      E="  -1.96%  "},
    @{Row=30; D="0.05356";    E="  -4.02%  "},
    @{Row=31; D="1.290";      E="  -1.19%  "},
    @{Row=32; D="3.469";      E="  -3.06%  "},
    @{Row=33; D="3.437";      E="  -0.36%  "},
    @{Row=34; D=$null;        E="  -1.55%  "},
    @{Row=35; D="2.875";      E=$null},
    @{Row=36; D="0.9520";     E="  -1.75%  "},
    @{Row=37; D="2.395";      E="  -1.04%  "},
    @{Row=38; D="0.5841";     E="  -2.25%  "},
    @{Row=39; D="0.01634";    E="  -1.18%  "},
    @{Row=40; D="1.112.38";   E="  +5.25%  "},
    @{Row=41; D="5.798";      E="  -2.06%  "},
    @{Row=42; D=$null;        E="  +0.03%  "},
    @{Row=43; D="0.8398";     E="  -1.65%  "},
    @{Row=44; D="101.03";     E="  -0.48%  "},
    @{Row=45; D="1.853.90";   E="  -0.76%  "},
    @{Row=46; D=$null;        E="  +0.47%  "},
    @{Row=47; D="57.74";      E="  -2.25%  "},
    @{Row=48; D="0.4542";     E="  +2.37%  "},
    @{Row=49; D="1.004";      E="  -0.12%  "},
    @{Row=50; D="8.134";      E="  -1.54%  "},
    @{Row=51; D="0.05226";    E="  -0.49%  "}
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        # Force these as plain text (many look like numbers, e.g. "1.005"),
        # matching the original inline-string cell type, then restore the
        # cell's style so no stray number-format style is left behind.
        $cell = $ws.Cells.Item($u.Row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}
